# Append the new 2022-05-03 (serial 44684) slate of games to Sheet1.
# Mirrors the existing table layout: date | visitor | home | visitor_odds | home_odds

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 261

$newRows = @(
    @{ Row = 262; Date = 44684; B = "Atlanta Braves";          C = "New York Mets";           D = $null; E = $null },
    @{ Row = 263; Date = 44684; B = "Atlanta Braves";          C = "New York Mets";           D = $null; E = $null },
    @{ Row = 264; Date = 44684; B = "Arizona Diamondbacks";    C = "Miami Marlins";           D = $null; E = $null },
    @{ Row = 265; Date = 44684; B = "Cinncinatti Reds";        C = "Milwaukee Brewers";       D = -160;  E = 120   },
    @{ Row = 266; Date = 44684; B = "Washington Nationals";    C = "Colorado Rockies";        D = -125;  E = -115  },
    @{ Row = 267; Date = 44684; B = "San Francisco Giants";    C = "Los Angeles Dodgers";     D = -120;  E = -110  },
    @{ Row = 268; Date = 44684; B = "Minnesota Twins";         C = "Baltimore Orioles";       D = 135;   E = -180  },
    @{ Row = 269; Date = 44684; B = "New York Yankees";        C = "Toronto Blue Jays";       D = -110;  E = -120  },
    @{ Row = 270; Date = 44684; B = "Los Angeles Angels";      C = "Boston Red Sox";          D = 105;   E = -140  },
    @{ Row = 271; Date = 44684; B = "Seattle Mariners";        C = "Houston Astros";          D = -125;  E = -105  },
    @{ Row = 272; Date = 44684; B = "Tampa Bay Rays";          C = "Oakland Athletics";       D = 115;   E = -155  },
    @{ Row = 273; Date = 44684; B = "San Diego Padres";        C = "Cleveland Gaurdians";     D = $null; E = $null },
    @{ Row = 274; Date = 44684; B = "Texas Rangers";           C = "Philidelphia Phillies";   D = -130;  E = -105  },
    @{ Row = 275; Date = 44684; B = "Pittsburgh Pirates";      C = "Detroit Tigers";          D = -125;  E = -110  },
    @{ Row = 276; Date = 44684; B = "Chicago White Sox";       C = "Chicago Cubs";            D = 115;   E = -155  },
    @{ Row = 277; Date = 44684; B = "Saint Luis Cardinals";    C = "Kansas City Royals";      D = 105;   E = -140  }
)

foreach ($r in $newRows) {
    # Clone formatting (incl. the date number format on column A) from the
    # preceding row so the new rows inherit the same style indices instead
    # of minting new ones.
    $ws.Range("A$lastRow`:E$lastRow").Copy($ws.Range("A$($r.Row):E$($r.Row)"))

    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C

    if ($null -eq $r.D) {
        $ws.Cells.Item($r.Row, 4).ClearContents()
    } else {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    }

    if ($null -eq $r.E) {
        $ws.Cells.Item($r.Row, 5).ClearContents()
    } else {
        $ws.Cells.Item($r.Row, 5).Value = $r.E
    }
}

# Match the author's final scroll position / selection from the diff.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 272
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G276").Select()
